# Update "想去人数" (want-to-go count) figures scraped for gh-pages output.
# Mirrors the refreshed crawl results across the 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 8600
$ws1.Range("F8").Value = 222
$ws1.Range("F12").Value = 5829
$ws1.Range("F15").Value = 373
$ws1.Range("F18").Value = 263
$ws1.Range("F23").Value = 9660
$ws1.Range("F25").Value = 1787
$ws1.Range("F27").Value = 0
$ws1.Range("F37").Value = 289
$ws1.Range("F38").Value = 1293
$ws1.Range("F45").Value = 164
$ws1.Range("F47").Value = 1060

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 12
$ws2.Range("F6").Value = 0

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 12
$ws4.Range("F8").Value = 8600
$ws4.Range("F10").Value = 222
$ws4.Range("F16").Value = 5829
$ws4.Range("F17").Value = 5829
$ws4.Range("F20").Value = 373
$ws4.Range("F22").Value = 263
$ws4.Range("F27").Value = 9660
$ws4.Range("F29").Value = 1787
$ws4.Range("F30").Value = 1468
$ws4.Range("F38").Value = 289
$ws4.Range("F39").Value = 1293
$ws4.Range("F46").Value = 164
$ws4.Range("F48").Value = 1060
